$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff adds "NA" values in column E (duplicate_image_filename) for
# all the data rows (2 through 21) that previously had no value there.
$ws.Range("E2:E21").Value = "NA"
